$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "level": update the randomly-generated level-design grid (cols A-D)
# while leaving the formula column E untouched.
# ---------------------------------------------------------------------------
$level = $wb.Worksheets.Item("level")

# Row 1
$level.Range("C1").ClearContents()
$level.Range("D1").Value = 8

# Row 2
$level.Range("B2").ClearContents()

# Row 4
$level.Range("A4").ClearContents()
$level.Range("D4").ClearContents()
$level.Range("B4").Value = 7
$level.Range("C4").Value = 7

# Row 5
$level.Range("C5").ClearContents()

# Row 7
$level.Range("A7").Value = 9
$level.Range("C7").Value = 8
$level.Range("D7").Value = 9

# Row 8
$level.Range("C8").ClearContents()

# Row 10
$level.Range("A10").Value = 7
$level.Range("B10").Value = 8
$level.Range("C10").Value = 9

# Row 12
$level.Range("A12").ClearContents()
$level.Range("C12").ClearContents()

# Row 14
$level.Range("D14").ClearContents()

# Row 15
$level.Range("A15").ClearContents()
$level.Range("B15").ClearContents()

# Row 17
$level.Range("C17").ClearContents()

# Row 18
$level.Range("A18").ClearContents()
$level.Range("B18").ClearContents()
$level.Range("D18").ClearContents()

# Row 19
$level.Range("C19").ClearContents()

# Row 20
$level.Range("A20").ClearContents()

# Row 21
$level.Range("D21").ClearContents()

# ---------------------------------------------------------------------------
# Sheet "enemies": rebalance weights/speeds, add two new resource-path rows
# and a new "all enemies" row (9).
# ---------------------------------------------------------------------------
$enemies = $wb.Worksheets.Item("enemies")

# Row 2 (index 1)
$enemies.Range("C2").Value = 0
$enemies.Range("D2").Value = 2
$enemies.Range("F2").Value = "EnemyPrefabs/Special Enemies/Halloween Bee/Halloween Bee"

# Row 3 (index 2)
$enemies.Range("E3").Value = 1
$enemies.Range("F3").Value = "EnemyPrefabs/Bullet Enemies/Neo Fly/Neo Fly"

# Row 4 (index 3)
$enemies.Range("C4").Value = 0
$enemies.Range("D4").Value = 2
$enemies.Range("F4").Value = "EnemyPrefabs/Special Enemies/Bionic Lady Bird/Bionic Lady Bird"

# Row 5 (index 4)
$enemies.Range("E5").Value = 1.5

# Row 6 (index 5)
$enemies.Range("B6").Value = 0
$enemies.Range("C6").Value = 0
$enemies.Range("D6").Value = 3
$enemies.Range("F6").Value = "EnemyPrefabs/Bullet Enemies/Neo Fly/Neo Fly"
$enemies.Range("G6").Value = "EnemyPrefabs/Bullet Enemies//"

# Row 7 (index 6)
$enemies.Range("D7").Value = 3
$enemies.Range("E7").Value = 1.5
$enemies.Range("G7").Value = "EnemyPrefabs/Special Enemies//"

# Row 8 (index 7)
$enemies.Range("B8").Value = 4
$enemies.Range("E8").Value = 1
$enemies.Range("F8").Value = "EnemyPrefabs/Arrow Enemies/Bee/Bee Arrow"

# Row 9 (index 8)
$enemies.Range("B9").Value = 0
$enemies.Range("C9").Value = 2
$enemies.Range("E9").Value = 1
$enemies.Range("F9").Value = "EnemyPrefabs/Bullet Enemies/Neo Fly/Neo Fly"

# Row 10 (index 9, brand new)
$enemies.Range("A10").Value = 9
$enemies.Range("B10").Value = 0
$enemies.Range("C10").Value = 0
$enemies.Range("D10").Value = 2
$enemies.Range("E10").Value = 1
$enemies.Range("F10").Value = "EnemyPrefabs/Special Enemies/Steampunk Fly/Steampunk Fly"

# ---------------------------------------------------------------------------
# Sheet "misc": fill in the previously-empty row 5 of the "all enemies" list.
# ---------------------------------------------------------------------------
$misc = $wb.Worksheets.Item("misc")
$misc.Range("A5").Value = "EnemyPrefabs/Arrow Enemies/Butterfly/Butterfly Arrow"

# ---------------------------------------------------------------------------
# Restore the per-sheet selections (cosmetic, matches the saved workbook
# state). "meta" is selected last so it stays the active/visible tab.
# ---------------------------------------------------------------------------
$level.Range("E14").Select()
$enemies.Range("F12").Select()
$misc.Range("A49").Select()
$meta = $wb.Worksheets.Item("meta")
$meta.Range("E5").Select()

Write-Output "edit applied"
